# Add month/year date formats to the split-workbook dates sample.
# Rows 8-13 already exist (B column has values formatted with various date
# formats) - we add the matching format-name labels in column A.
# Rows 14-19 are brand new: more date-format examples plus a couple of
# empty placeholder rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Label the existing rows 8-13 with their format names (column A) ---
# (Row 12's label is entered before row 11's so new shared-strings land in
# the same order as the original authoring session.)
$ws.Cells.Item(8, 1).Value = "D/MM/YY"
$ws.Cells.Item(9, 1).Value = "D/M/YY"
$ws.Cells.Item(10, 1).Value = "DD/MM/YY"
$ws.Cells.Item(12, 1).Value = "DD-MONTH-YYYY"
$ws.Cells.Item(11, 1).Value = "DD-MON-YY"
$ws.Cells.Item(13, 1).Value = "D Month YYYY"

# --- New row 14: MM/DD/YY example ---
$ws.Cells.Item(14, 1).Value = "MM/DD/YY"
$ws.Cells.Item(14, 2).Value = 32
$ws.Cells.Item(14, 2).NumberFormat = "mm/dd/yy;@"

# --- New row 15: ISO yyyy-mm-dd example (no label) ---
$ws.Cells.Item(15, 2).Value = 32
$ws.Cells.Item(15, 2).NumberFormat = "yyyy-mm-dd;@"

# --- New row 16: yy/mm/dd example (no label) ---
$ws.Cells.Item(16, 2).Value = 32
$ws.Cells.Item(16, 2).NumberFormat = "yy/mm/dd;@"

# --- New row 17: reuse the existing yyyy/mm/dd format ---
$ws.Cells.Item(17, 2).Value = 32
$ws.Cells.Item(17, 2).NumberFormat = "yyyy/mm/dd;@"

# --- New rows 18-19: empty cells using the existing m/d/yyyy format ---
$ws.Cells.Item(18, 2).NumberFormat = "m/d/yyyy;@"
$ws.Cells.Item(19, 2).NumberFormat = "m/d/yyyy;@"

# --- Widen column A so the new, longer labels are visible ---
$ws.Columns.Item(1).ColumnWidth = 15.592447916666666

# --- Move the active selection the way it ended up after the edits ---
$ws.Range("D21").Select()
